# The paragraph contains an orange (FF6600) run reading
#   " Remember to update xls file. "
# The commit splits it (dropping the trailing space) and appends a new
# green (008000) note plus some extra orange spaces, i.e. the run content
# becomes the concatenation of:
#   " Remember to update xls file."   (orange)
#   " "                               (orange)
#   "[ _v15 has been updated. ]"      (green)
#   " "                               (orange)
#   " "                               (orange)

$d = $word.ActiveDocument

$wdColorOrange = 26367   # RGB(0xFF,0x66,0x00) packed as BGR (wdColor value)
$wdColorGreen  = 32768   # RGB(0x00,0x80,0x00) packed as BGR (wdColor value)

$rng = $d.Content
$found = $rng.Find.Execute(" Remember to update xls file. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Overwrite the matched text with the full new text; the new range
    # inherits the matched run's (orange) formatting.
    $newText = " Remember to update xls file. [ _v15 has been updated. ]  "
    $rng.Text = $newText
    $baseStart = $rng.Start

    # Segment lengths within $newText.
    $seg1 = " Remember to update xls file."  # orange
    $seg2 = " "                              # orange
    $seg3 = "[ _v15 has been updated. ]"     # green
    $seg4 = " "                              # orange
    $seg5 = " "                              # orange

    $o1 = 0
    $o2 = $o1 + $seg1.Length
    $o3 = $o2 + $seg2.Length
    $o4 = $o3 + $seg3.Length
    $o5 = $o4 + $seg4.Length
    $o6 = $o5 + $seg5.Length

    # Color each segment explicitly, in order, so the run boundaries line up
    # with the segments above.
    $d.Range($baseStart + $o1, $baseStart + $o2).Font.Color = $wdColorOrange
    $d.Range($baseStart + $o2, $baseStart + $o3).Font.Color = $wdColorOrange
    $d.Range($baseStart + $o3, $baseStart + $o4).Font.Color = $wdColorGreen
    $d.Range($baseStart + $o4, $baseStart + $o5).Font.Color = $wdColorOrange
    $d.Range($baseStart + $o5, $baseStart + $o6).Font.Color = $wdColorOrange
}
